# Apply the described changes to the workbook.
$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# Sheet "Ingredients": swap columns A (id) and B (name)
# ----------------------------------------------------------------------
$wsIngredients = $wb.Worksheets.Item("Ingredients")

# Header row: swap A1/B1 text
$wsIngredients.Range("A1").Value = "name"
$wsIngredients.Range("B1").Value = "id"

# Data rows 2..11: swap id/name values (A<->B) for each row
$ingredientRows = @(
    @{Id = 9;  Name = "Ice Cream"},
    @{Id = 10; Name = "Yogurt"},
    @{Id = 11; Name = "Bread"},
    @{Id = 12; Name = "Oats"},
    @{Id = 13; Name = "Rice"},
    @{Id = 14; Name = "Boiled Potatoes"},
    @{Id = 15; Name = "Red Chili"},
    @{Id = 16; Name = "Milk"},
    @{Id = 17; Name = "Flour"},
    @{Id = 21; Name = "Tomatoes"}
)

$r = 2
foreach ($row in $ingredientRows) {
    $wsIngredients.Cells.Item($r, 1).Value = $row.Name
    $wsIngredients.Cells.Item($r, 2).Value = $row.Id
    $r = $r + 1
}

# ----------------------------------------------------------------------
# Sheet "Recipes": reorder columns to name/steps/id/ingredients,
# and replace the 3 existing data rows with a single new recipe row.
# ----------------------------------------------------------------------
$wsRecipes = $wb.Worksheets.Item("Recipes")

# Clear out old data rows (2..4) first
$wsRecipes.Range("A2:D4").Clear()

# Header row
$wsRecipes.Range("A1").Value = "name"
$wsRecipes.Range("B1").Value = "steps"
$wsRecipes.Range("C1").Value = "id"
$wsRecipes.Range("D1").Value = "ingredients"

# New single data row
$steps = "Measure 1 cup of sugar into a mixing bowl.`nMeasure 1 cup of butter into the same mixing bowl.`nMix the sugar and butter together until combined.`nAdd eggs to the same bowl and mix until combined.`nMeasure 2 teaspoons of vanilla extract into the same mixing bowl.`nMix the vanilla extract until combined. `nMeasure 1.5 cups of all purpose flour into a new mixing bowl.`nMeasure 1.75 teaspoons of baking powder into bowl with the flour.`nBake.`n"
$ingredients = "1 cup white sugar, 0.5 cups unsalted butter, 2 teaspoons vanilla extract, 1.5 cups all purpose flour, 1.75 teaspoons baking powder, 0.5 cups milk"

$wsRecipes.Range("A2").Value = "White Cake"
$wsRecipes.Range("B2").Value = $steps
$wsRecipes.Range("C2").Value = 20
$wsRecipes.Range("D2").Value = $ingredients
